# Actualización automática de tasas-transfi.xlsx
# Updates the "tasas" sheet rate cells (N10, O10, N12, O12) and refreshes the
# manually-summarized daily conversion text on "Hoja1" (cell A1) to match.

$wb = $excel.ActiveWorkbook

$wsTasas = $wb.Worksheets.Item("tasas")
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

# New rate values for the "tasas" sheet.
$wsTasas.Range("N10").Value = 80.25
$wsTasas.Range("O10").Value = 4093.11
$wsTasas.Range("N12").Value = 4085
$wsTasas.Range("O12").Value = 78.09999999999999

# Recompute the Binance summary lines shown in the Hoja1!A1 note, following
# the same rounding used by the previous manual update. Intermediate
# (unrounded) results chain into the next step; only the displayed text is
# rounded to 2 decimals:
#   rate1raw  = 1000 / N10          -> displayed as rate1
#   pesos1raw = rate1raw * O10      -> displayed as pesos1
#   rate2raw  = pesos1raw / N12     -> displayed as rate2
#   bs2raw    = O12 * pesos1raw / N12 -> displayed as bs2
$rate1raw  = 1000 / 80.25
$pesos1raw = $rate1raw * 4093.11
$rate2raw  = $pesos1raw / 4085
$bs2raw    = 78.09999999999999 * $pesos1raw / 4085

$rate1  = [Math]::Round($rate1raw, 2)
$pesos1 = [Math]::Round($pesos1raw, 2)
$rate2  = [Math]::Round($rate2raw, 2)
$bs2    = [Math]::Round($bs2raw, 2)

$oldText = $wsHoja1.Range("A1").Value()

# Match the two "Binance" summary lines regardless of the exact numbers they
# currently hold, then substitute in the freshly computed ones.
$pattern = "✅ 1000 Bs = [\d.,]+ = [\d.,]+ pesos`n✅ [\d.,]+ pesos = [\d.,]+ = [\d.,]+ Bs"
$replacement = "✅ 1000 Bs = $rate1 = $pesos1 pesos`n✅ $pesos1 pesos = $rate2 = $bs2 Bs"
$newText = [regex]::Replace($oldText, $pattern, $replacement)

$wsHoja1.Range("A1").Value = $newText
